$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44 (shifts existing rows 44..95 down to 45..96)
$ws.Rows(44).Insert()

# Populate the newly inserted row 44 with the new weekly price record
$ws.Range("A44").Value = 11
$ws.Range("B44").Value = "Vega Monumental Concepción"
$ws.Range("C44").Value = "Bíobío"
$ws.Range("D44").Value = 44981
$ws.Range("E44").Value = 8
$ws.Range("F44").Value = 100112037
$ws.Range("G44").Value = "Cebollín"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 100
$ws.Range("K44").Value = 6000
$ws.Range("L44").Value = 6500
$ws.Range("M44").Value = 6250
$ws.Range("N44").Value = "$/paquete 36 unidades"
$ws.Range("O44").Value = "Región Metropolitana"
$ws.Range("P44").Value = 174
$ws.Range("Q44").Value = 36
$ws.Range("R44").Value = "Hortaliza"
